$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.125.83'
$ws.Range("E2").Value = '  +0.68%  '

$ws.Range("D3").Value = '3.850.74'
$ws.Range("E3").Value = '  +1.11%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '691.10'
$ws.Range("E5").Value = '  +3.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.30'
$ws.Range("E6").Value = '  +2.19%  '

$ws.Range("D7").Value = '3.848.34'
$ws.Range("E7").Value = '  +1.06%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  +0.13%  '

$ws.Range("E10").Value = '  +1.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.40'
$ws.Range("E11").Value = '  +5.74%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.461'
$ws.Range("E12").Value = '  -0.49%  '

$ws.Range("E13").Value = '  +5.68%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.65'
$ws.Range("E14").Value = '  +1.99%  '

$ws.Range("D15").Value = '4.496.09'
$ws.Range("E15").Value = '  +0.84%  '

$ws.Range("D16").Value = '3.850.63'
$ws.Range("E16").Value = '  +0.96%  '

$ws.Range("D17").Value = '71.148.35'
$ws.Range("E17").Value = '  +0.70%  '

$ws.Range("E18").Value = '  +0.48%  '

$ws.Range("E19").Value = '  +0.50%  '

$ws.Range("E20").Value = '  +0.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.16'
$ws.Range("E21").Value = '  -4.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '488.72'
$ws.Range("E22").Value = '  +2.87%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.721'
$ws.Range("E23").Value = '  +0.68%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.83'
$ws.Range("E24").Value = '  +2.04%  '

$ws.Range("E25").Value = '  +1.95%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.38'
$ws.Range("E26").Value = '  +1.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.50'
$ws.Range("E27").Value = '  +1.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.15'
$ws.Range("E28").Value = '  +0.73%  '

$ws.Range("D29").Value = '4.001.84'
$ws.Range("E29").Value = '  +0.93%  '

$ws.Range("E30").Value = '  +0.10%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.11'
$ws.Range("E31").Value = '  +8.43%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.33'
$ws.Range("E32").Value = '  +0.77%  '

$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.64'
$ws.Range("E33").Value = '  +3.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.76'
$ws.Range("E34").Value = '  +0.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.182'
$ws.Range("E35").Value = '  +3.62%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.27'
$ws.Range("E36").Value = '  +1.24%  '

$ws.Range("D37").Value = '3.799.45'
$ws.Range("E37").Value = '  +0.78%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.16%  '

$ws.Range("E39").Value = '  +1.32%  '

$ws.Range("E40").Value = '  +13.24%  '

$ws.Range("E41").Value = '  +0.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.05'
$ws.Range("E42").Value = '  +0.49%  '

$ws.Range("E43").Value = '  +4.61%  '

$ws.Range("E44").Value = '  -0.10%  '

$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '164.98'
$ws.Range("E46").Value = '  +3.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000307'
$ws.Range("E47").Value = '  +6.32%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.71'
$ws.Range("E48").Value = '  +1.57%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.46'
$ws.Range("E49").Value = '  -2.81%  '

$ws.Range("E50").Value = '  +0.96%  '

$ws.Range("E51").Value = '  -2.80%  '
